# Update loading_percent values on Sheet1 for the 380 kV case run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9.33390393874344
$ws.Range("C2").Value = 5.405776359788273
$ws.Range("D2").Value = 5.987150653745522
$ws.Range("E2").Value = 16.423163734356
$ws.Range("G2").Value = 3.639287268732075
$ws.Range("I2").Value = 21.62300519657893
$ws.Range("K2").Value = 9.379349347190169
$ws.Range("O2").Value = 23.37550271558159
$ws.Range("B3").Value = 9.004824561083781
$ws.Range("C3").Value = 5.096452989557434
$ws.Range("D3").Value = 5.868876234848567
$ws.Range("E3").Value = 15.49653435862805
$ws.Range("G3").Value = 3.641525777724216
$ws.Range("I3").Value = 21.7011353661785
$ws.Range("K3").Value = 9.145923172151665
$ws.Range("O3").Value = 23.42119823758021
$ws.Range("B4").Value = 8.798175871109006
$ws.Range("C4").Value = 4.895650809597979
$ws.Range("D4").Value = 5.796833810510158
$ws.Range("E4").Value = 14.90328505129345
$ws.Range("G4").Value = 3.642971211485527
$ws.Range("I4").Value = 21.75393237056602
$ws.Range("K4").Value = 9.00130500290066
$ws.Range("O4").Value = 23.45523815885859
$ws.Range("B5").Value = 8.712946002104598
$ws.Range("C5").Value = 4.811110159967606
$ws.Range("D5").Value = 5.767666612929891
$ws.Range("E5").Value = 14.65569392595479
$ws.Range("G5").Value = 3.643578148113814
$ws.Range("I5").Value = 21.77665673415794
$ws.Range("K5").Value = 8.942136299448443
$ws.Range("O5").Value = 23.47060819394613
$ws.Range("B6").Value = 8.698736329155455
$ws.Range("C6").Value = 4.796909234945989
$ws.Range("D6").Value = 5.762836266233636
$ws.Range("E6").Value = 14.61423804735643
$ws.Range("G6").Value = 3.643680013061414
$ws.Range("I6").Value = 21.78050300745034
$ws.Range("K6").Value = 8.932299846179829
$ws.Range("O6").Value = 23.47325070012559
$ws.Range("B7").Value = 8.7970303677497
$ws.Range("C7").Value = 4.894521611415921
$ws.Range("D7").Value = 5.79643961943814
$ws.Range("E7").Value = 14.89996918717937
$ws.Range("G7").Value = 3.642979324244022
$ws.Range("I7").Value = 21.75423394893872
$ws.Range("K7").Value = 9.000507865800049
$ws.Range("O7").Value = 23.45543938542062
$ws.Range("B8").Value = 9.22146832071858
$ws.Range("C8").Value = 5.301390886209933
$ws.Range("D8").Value = 5.946276522719229
$ws.Range("E8").Value = 16.1088519212923
$ws.Range("G8").Value = 3.640044412140663
$ws.Range("I8").Value = 21.64894066747799
$ws.Range("K8").Value = 9.299185010884623
$ws.Range("O8").Value = 23.39001340887794
$ws.Range("B9").Value = 10.01200329271725
$ws.Range("C9").Value = 6.012246178392517
$ws.Range("D9").Value = 6.242743042320279
$ws.Range("E9").Value = 18.33882573495292
$ws.Range("G9").Value = 3.63484941595157
$ws.Range("I9").Value = 21.48092341388828
$ws.Range("K9").Value = 9.870903238656476
$ws.Range("O9").Value = 23.30942822428655
$ws.Range("B10").Value = 10.5610328859374
$ws.Range("C10").Value = 6.480831297276842
$ws.Range("D10").Value = 6.459624067527034
$ws.Range("E10").Value = 19.96993774225593
$ws.Range("G10").Value = 3.631370280731129
$ws.Range("I10").Value = 21.38117277012684
$ws.Range("K10").Value = 10.27777012421788
$ws.Range("O10").Value = 23.27961700287941
$ws.Range("B11").Value = 10.80280074348471
$ws.Range("C11").Value = 6.682283179703724
$ws.Range("D11").Value = 6.557578797753612
$ws.Range("E11").Value = 20.67010341077044
$ws.Range("G11").Value = 3.629859995515184
$ws.Range("I11").Value = 21.34098913140645
$ws.Range("K11").Value = 10.45912440572677
$ws.Range("O11").Value = 23.27249258161217
$ws.Range("B12").Value = 10.89312584626401
$ws.Range("C12").Value = 6.756881514079188
$ws.Range("D12").Value = 6.59453059096987
$ws.Range("E12").Value = 20.92925844981496
$ws.Range("G12").Value = 3.629298434790165
$ws.Range("I12").Value = 21.32652356250381
$ws.Range("K12").Value = 10.52719604788787
$ws.Range("O12").Value = 23.27072385875772
$ws.Range("B13").Value = 10.87372844218456
$ws.Range("C13").Value = 6.740890493842482
$ws.Range("D13").Value = 6.58657927074365
$ws.Range("E13").Value = 20.87371031973763
$ws.Range("G13").Value = 3.62941891747711
$ws.Range("I13").Value = 21.32960550859907
$ws.Range("K13").Value = 10.51256343704799
$ws.Range("O13").Value = 23.27106341134107
$ws.Range("B14").Value = 10.81025685236574
$ws.Range("C14").Value = 6.688454253515895
$ws.Range("D14").Value = 6.560621902375599
$ws.Range("E14").Value = 20.6915439597407
$ws.Range("G14").Value = 3.629813588448289
$ws.Range("I14").Value = 21.33978396283484
$ws.Range("K14").Value = 10.46473715408512
$ws.Range("O14").Value = 23.27232842830271
$ws.Range("B15").Value = 10.77121670035048
$ws.Range("C15").Value = 6.656115767092237
$ws.Range("D15").Value = 6.544702680422702
$ws.Range("E15").Value = 20.57918369226482
$ws.Range("G15").Value = 3.630056682164799
$ws.Range("I15").Value = 21.34611650171636
$ws.Range("K15").Value = 10.4353617227798
$ws.Range("O15").Value = 23.27322438365071
$ws.Range("B16").Value = 10.54506511278681
$ws.Range("C16").Value = 6.467429984893864
$ws.Range("D16").Value = 6.45320482221738
$ws.Range("E16").Value = 19.92334129240778
$ws.Range("G16").Value = 3.631470433127354
$ws.Range("I16").Value = 21.38390373948419
$ws.Range("K16").Value = 10.26583728077716
$ws.Range("O16").Value = 23.28021244200033
$ws.Range("B17").Value = 10.4042234690526
$ws.Range("C17").Value = 6.348675204924064
$ws.Range("D17").Value = 6.396863939665239
$ws.Range("E17").Value = 19.51031199713599
$ws.Range("G17").Value = 3.632356223066126
$ws.Range("I17").Value = 21.40841822956597
$ws.Range("K17").Value = 10.16083459589267
$ws.Range("O17").Value = 23.28615072453158
$ws.Range("B18").Value = 10.32246720278073
$ws.Range("C18").Value = 6.279268458307828
$ws.Range("D18").Value = 6.36439401018348
$ws.Range("E18").Value = 19.26880874776455
$ws.Range("G18").Value = 3.632872523235681
$ws.Range("I18").Value = 21.42300689409556
$ws.Range("K18").Value = 10.10009223989855
$ws.Range("O18").Value = 23.29017196113051
$ws.Range("B19").Value = 10.29466003784667
$ws.Range("C19").Value = 6.255579366921592
$ws.Range("D19").Value = 6.353390490067025
$ws.Range("E19").Value = 19.18636222529449
$ws.Range("G19").Value = 3.633048506246625
$ws.Range("I19").Value = 21.42803014688818
$ws.Range("K19").Value = 10.07946834554324
$ws.Range("O19").Value = 23.2916373938648
$ws.Range("B20").Value = 10.41929430402675
$ws.Range("C20").Value = 6.361430951259379
$ws.Range("D20").Value = 6.402868458765782
$ws.Range("E20").Value = 19.55468717491102
$ws.Range("G20").Value = 3.63226122407948
$ws.Range("I20").Value = 21.40575802279793
$ws.Range("K20").Value = 10.17204878566455
$ws.Range("O20").Value = 23.28545586964485
$ws.Range("B21").Value = 10.82893387783104
$ws.Range("C21").Value = 6.70390185236198
$ws.Range("D21").Value = 6.568250360451788
$ws.Range("E21").Value = 20.74521273404633
$ws.Range("G21").Value = 3.629697383542351
$ws.Range("I21").Value = 21.3367738855739
$ws.Range("K21").Value = 10.47880177021342
$ws.Range("O21").Value = 23.27193162042439
$ws.Range("B22").Value = 11.08947462592464
$ws.Range("C22").Value = 6.917892382433813
$ws.Range("D22").Value = 6.67549499734226
$ws.Range("E22").Value = 21.48843472504073
$ws.Range("G22").Value = 3.628082075276368
$ws.Range("I22").Value = 21.29606869771921
$ws.Range("K22").Value = 10.67573905725649
$ws.Range("O22").Value = 23.26850969040441
$ws.Range("B23").Value = 10.95109974917724
$ws.Range("C23").Value = 6.804581993978454
$ws.Range("D23").Value = 6.618346224868477
$ws.Range("E23").Value = 21.09494145361025
$ws.Range("G23").Value = 3.628938696501899
$ws.Range("I23").Value = 21.31739166206512
$ws.Range("K23").Value = 10.57097477861135
$ws.Range("O23").Value = 23.2698393797072
$ws.Range("B24").Value = 10.41248321709989
$ws.Range("C24").Value = 6.355667607335912
$ws.Range("D24").Value = 6.400154058094573
$ws.Range("E24").Value = 19.5346377664331
$ws.Range("G24").Value = 3.632304151162149
$ws.Range("I24").Value = 21.40695916073091
$ws.Range("K24").Value = 10.16698001308811
$ws.Range("O24").Value = 23.28576812217514
$ws.Range("B25").Value = 9.803312853516003
$ws.Range("C25").Value = 5.829336501748694
$ws.Range("D25").Value = 6.162526492926562
$ws.Range("E25").Value = 17.71340042964233
$ws.Range("G25").Value = 3.636195222637485
$ws.Range("I25").Value = 21.52223343119667
$ws.Range("K25").Value = 9.718235449876882
$ws.Range("O25").Value = 23.32608670815555
